# Commit: "Fruta / hortaliza, semanal"
# A new weekly price-report row is inserted into the daily logic sheet at
# row 1151, pushing the former rows 1151-1215 down to 1152-1216 (the sheet's
# used range grows from A1:R1215 to A1:R1216).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 1151; Excel shifts rows 1151:1215 down to 1152:1216
# and carries the date column's number format (style) along with it.
$ws.Rows.Item(1151).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A1151").Value = 6
$ws.Range("B1151").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C1151").Value = "Metropolitana"
$ws.Range("D1151").Value = "2022-09-22"
$ws.Range("E1151").Value = 13
$ws.Range("F1151").Value = 100112024
$ws.Range("G1151").Value = "Choclo"
$ws.Range("H1151").Value = "Dulce o Americano"
$ws.Range("I1151").Value = "Primera"
$ws.Range("J1151").Value = 460
$ws.Range("K1151").Value = 28000
$ws.Range("L1151").Value = 28000
$ws.Range("M1151").Value = 28000
$ws.Range("N1151").Value = "$/malla 70 unidades"
$ws.Range("O1151").Value = "Región de Arica y Parinacota"
$ws.Range("P1151").Value = 400
$ws.Range("Q1151").Value = 70
$ws.Range("R1151").Value = "Hortaliza"
